$d = $word.ActiveDocument

# 1. Replace the title / heading text everywhere it appears (Heading1 title
#    and the bold repeated title near the end of the doc).
$d.Content.Find.Execute(
    "Play Bingo Billions Free: Fun Gameplay and Great Winning Potential",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Bingo Billions Free - Exciting Gameplay and Big Wins", 2) | Out-Null

# 2. Replace the italic meta-description paragraph text.
$d.Content.Find.Execute(
    "Read our review of Bingo Billions, a bingo-themed slot game with fun gameplay and great winning potential. Play free and enjoy frequent payouts and free spins.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Read our review of Bingo Billions and play for free. Exciting gameplay features and big winning potential.", 2) | Out-Null

# 3. Update the first bullet under "What we like". Rebuild the whole
#    paragraph via InsertXML (instead of Find/Replace) so the existing
#    leading empty <w:r/> run stays intact rather than being collapsed.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Fun gameplay features\r?$") {
        $xmlFun = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Exciting gameplay features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xmlFun) | Out-Null
        break
    }
}

# 4. Remove the "Potential for great wins" and "Frequent payouts in base
#    game" bullets from the "What we like" list.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "^Potential for great wins\r?$" -or $t -match "^Frequent payouts in base game\r?$") {
        $p.Range.Delete()
    }
}

# 5. Append two new bullets after "Free Spins feature with tripled prizes":
#    "Wild and Scatter symbols with payouts" and
#    "Three bonus symbols for added excitement".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Free Spins feature with tripled prizes\r?$") {
        $p.Range.InsertParagraphAfter() | Out-Null
        $p1 = $d.Paragraphs.Item($i + 1)
        $xml1 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Wild and Scatter symbols with payouts</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p1.Range.InsertXML($xml1) | Out-Null

        $p1 = $d.Paragraphs.Item($i + 1)
        $p1.Range.InsertParagraphAfter() | Out-Null
        $p2 = $d.Paragraphs.Item($i + 2)
        $xml2 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Three bonus symbols for added excitement</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p2.Range.InsertXML($xml2) | Out-Null
        break
    }
}

Write-Output "done"
